# Auto-generated edit script applying the Kraken_Profits H:N value updates
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H9").Value = 418.75
$ws.Range("I9").Value = 516.6667
$ws.Range("K9").Value = 516.6667
$ws.Range("M9").Value = -347.6667

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H100").Value = 13000
$ws.Range("J100").Value = 13000
$ws.Range("L100").Value = 13000
$ws.Range("N100").Value = -14082

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H22").Value = 1011.75
$ws.Range("I22").Value = 349
$ws.Range("J22").Value = 3000
$ws.Range("K22").Value = 349
$ws.Range("L22").Value = 3000
$ws.Range("M22").Value = -50
$ws.Range("N22").Value = -3598

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H44").Value = 35000
$ws.Range("I44").Value = 0
$ws.Range("K44").Value = 0
$ws.Range("M44").ClearContents()

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 3500
$ws.Range("I61").Value = 3500
$ws.Range("K61").Value = 3500
$ws.Range("M61").Value = -3288

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H97").Value = 4160
$ws.Range("J97").Value = 4599.6665
$ws.Range("L97").Value = 4599.6665
$ws.Range("N97").Value = -5591.6665

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H132").Value = 4133.1665
$ws.Range("I132").Value = 2949.6667
$ws.Range("K132").Value = 8849.000100000001
$ws.Range("M132").Value = -6319.000100000001

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H136").Value = 3500
$ws.Range("I136").Value = 3500
$ws.Range("K136").Value = 10500
$ws.Range("M136").Value = -7950

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H99").Value = 2673.6667
$ws.Range("I99").Value = 3505
$ws.Range("J99").Value = 1011
$ws.Range("K99").Value = 3505
$ws.Range("L99").Value = 1011
$ws.Range("M99").Value = -2007
$ws.Range("N99").Value = -4007

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H105").Value = 4083.1667
$ws.Range("I105").Value = 3166.6667
$ws.Range("J105").Value = 4999.6665
$ws.Range("K105").Value = 3166.6667
$ws.Range("L105").Value = 4999.6665
$ws.Range("M105").Value = -1419.6667
$ws.Range("N105").Value = -8493.666499999999

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H134").Value = 9249
$ws.Range("I134").Value = 8665.666999999999
$ws.Range("K134").Value = 25997.001
$ws.Range("M134").Value = -23462.001

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H15").Value = 5432.6665
$ws.Range("I15").Value = 300
$ws.Range("J15").Value = 7999
$ws.Range("K15").Value = 300
$ws.Range("L15").Value = 7999
$ws.Range("M15").Value = -130
$ws.Range("N15").Value = -8339

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 7813.231
$ws.Range("J31").Value = 7267.7
$ws.Range("L31").Value = 7267.7
$ws.Range("N31").Value = -7857.7

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H34").Value = 7813.231
$ws.Range("J34").Value = 7267.7
$ws.Range("L34").Value = 7267.7
$ws.Range("N34").Value = -7671.7

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H36").Value = 5300
$ws.Range("I36").Value = 5300
$ws.Range("K36").Value = 5300
$ws.Range("M36").Value = -4912

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H40").Value = 5300
$ws.Range("I40").Value = 5300
$ws.Range("K40").Value = 5300
$ws.Range("M40").Value = -5140

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H132").Value = 1578.4286
$ws.Range("I132").Value = 1977.75
$ws.Range("J132").Value = 1046
$ws.Range("K132").Value = 5933.25
$ws.Range("L132").Value = 3138
$ws.Range("M132").Value = -3403.25
$ws.Range("N132").Value = -8198

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H2").Value = 50.333332
$ws.Range("I2").Value = 31.375
$ws.Range("K2").Value = 188.25
$ws.Range("M2").Value = -75.25

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H12").Value = 3822.2222
$ws.Range("I12").Value = 466
$ws.Range("K12").Value = 1398
$ws.Range("M12").Value = -1225

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H43").Value = 2500
$ws.Range("I43").Value = 1000
$ws.Range("J43").Value = 4000
$ws.Range("K43").Value = 3000
$ws.Range("L43").Value = 12000
$ws.Range("M43").Value = -2886
$ws.Range("N43").Value = -12228

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H95").Value = 80000
$ws.Range("I95").Value = 80000
$ws.Range("K95").Value = 240000
$ws.Range("M95").Value = -237941

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H13").Value = 1005
$ws.Range("I13").Value = 1005
$ws.Range("J13").Value = 0
$ws.Range("K13").Value = 1005
$ws.Range("L13").Value = 0
$ws.Range("M13").Value = -866
$ws.Range("N13").ClearContents()

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H57").Value = 19770
$ws.Range("J57").Value = 23448.75
$ws.Range("L57").Value = 23448.75
$ws.Range("N57").Value = -25088.75

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 4832.6665
$ws.Range("J70").Value = 4999.25
$ws.Range("L70").Value = 4999.25
$ws.Range("N70").Value = -5539.25

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H73").Value = 4832.6665
$ws.Range("J73").Value = 4999.25
$ws.Range("L73").Value = 4999.25
$ws.Range("N73").Value = -6871.25

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H123").Value = 39000
$ws.Range("J123").Value = 39000
$ws.Range("L123").Value = 39000
$ws.Range("N123").Value = -43900

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H132").Value = 5499.4
$ws.Range("I132").Value = 4624.5
$ws.Range("K132").Value = 13873.5
$ws.Range("M132").Value = -11343.5

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H16").Value = 1498.909
$ws.Range("I16").Value = 1148.8
$ws.Range("J16").Value = 5000
$ws.Range("K16").Value = 1148.8
$ws.Range("L16").Value = 5000
$ws.Range("M16").Value = -978.8
$ws.Range("N16").Value = -5340

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H46").Value = 3150.8
$ws.Range("J46").Value = 3150.8
$ws.Range("L46").Value = 3150.8
$ws.Range("N46").Value = -3526.8

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H68").Value = 1788.2307
$ws.Range("I68").Value = 1674.7
$ws.Range("J68").Value = 2166.6667
$ws.Range("K68").Value = 1674.7
$ws.Range("L68").Value = 2166.6667
$ws.Range("M68").Value = -925.7
$ws.Range("N68").Value = -3664.6667

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H71").Value = 1788.2307
$ws.Range("I71").Value = 1674.7
$ws.Range("J71").Value = 2166.6667
$ws.Range("K71").Value = 8373.5
$ws.Range("L71").Value = 10833.3335
$ws.Range("M71").Value = -4629.5
$ws.Range("N71").Value = -18321.3335

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H110").Value = 42500
$ws.Range("J110").Value = 42500
$ws.Range("L110").Value = 42500
$ws.Range("N110").Value = -50680

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H132").Value = 8000
$ws.Range("I132").Value = 8000
$ws.Range("K132").Value = 24000
$ws.Range("M132").Value = -21470

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H136").Value = 2470.5557
$ws.Range("I136").Value = 1583.75
$ws.Range("J136").Value = 3180
$ws.Range("K136").Value = 4751.25
$ws.Range("L136").Value = 9540
$ws.Range("M136").Value = -2201.25
$ws.Range("N136").Value = -14640
